$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers must be forced to
# remain text (matching the original inline-string cell type) by
# setting the number format to Text before assigning the value.

$ws.Range('D2').Value = '66.048.88'
$ws.Range('E2').Value = '  +6.71%  '
$ws.Range('D3').Value = '3.013.25'
$ws.Range('E3').Value = '  +3.98%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '585.33'
$ws.Range('E5').Value = '  +2.92%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '163.59'
$ws.Range('E6').Value = '  +13.62%  '
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').Value = '3.011.29'
$ws.Range('E8').Value = '  +3.98%  '
$ws.Range('E9').Value = '  +3.25%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.99'
$ws.Range('E10').Value = '  +1.34%  '
$ws.Range('E11').Value = '  +8.44%  '
$ws.Range('E12').Value = '  +6.49%  '
$ws.Range('E13').Value = '  +9.70%  '
$ws.Range('E14').Value = '  +8.05%  '
$ws.Range('E15').Value = '  +0.95%  '
$ws.Range('D16').Value = '66.036.68'
$ws.Range('E16').Value = '  +6.74%  '
$ws.Range('D17').Value = '3.512.68'
$ws.Range('E17').Value = '  +3.99%  '
$ws.Range('E18').Value = '  +6.93%  '
$ws.Range('D19').Value = '3.011.24'
$ws.Range('E19').Value = '  +4.08%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '458.78'
$ws.Range('E20').Value = '  +6.38%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.98'
$ws.Range('E21').Value = '  +7.76%  '
$ws.Range('E22').Value = '  +5.25%  '
$ws.Range('E23').Value = '  +7.65%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '82.50'
$ws.Range('E24').Value = '  +4.67%  '
$ws.Range('E25').Value = '  +14.22%  '
$ws.Range('E26').Value = '  +2.65%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.69'
$ws.Range('E27').Value = '  +6.24%  '
$ws.Range('E28').Value = '  -0.06%  '
$ws.Range('E29').Value = '  +15.86%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.34'
$ws.Range('E30').Value = '  +15.58%  '
$ws.Range('E31').Value = '  +4.53%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0000103'
$ws.Range('E32').Value = '  -6.07%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '27.02'
$ws.Range('E33').Value = '  +5.74%  '
$ws.Range('E34').Value = '  +3.24%  '
$ws.Range('E35').Value = '  -0.15%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.996'
$ws.Range('E36').Value = '  +4.48%  '
$ws.Range('E37').Value = '  +7.93%  '
$ws.Range('E38').Value = '  +11.91%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.07'
$ws.Range('E39').Value = '  +7.88%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '49.95'
$ws.Range('E40').Value = '  +2.32%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.307'
$ws.Range('E41').Value = '  +14.53%  '
$ws.Range('E42').Value = '  +6.46%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '43.84'
$ws.Range('E43').Value = '  +9.21%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.51'
$ws.Range('E44').Value = '  +4.61%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '388.25'
$ws.Range('E45').Value = '  +11.83%  '
$ws.Range('E46').Value = '  +6.76%  '
$ws.Range('D47').Value = '2.799.62'
$ws.Range('E47').Value = '  +3.80%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '135.37'
$ws.Range('E48').Value = '  +2.72%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '23.99'
$ws.Range('E50').Value = '  +10.81%  '
$ws.Range('E51').Value = '  +4.16%  '
